$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.571.05'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '1.670.53'
$ws.Range("E3").Value = '  +2.27%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''239.42'
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '''0.4792'
$ws.Range("E7").Value = '  +2.28%  '
$ws.Range("E8").Value = '  +3.42%  '
$ws.Range("D9").Value = '''0.06174'
$ws.Range("E9").Value = '  +3.28%  '
$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.669.95'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.06993'
$ws.Range("E11").Value = '  -2.59%  '
$ws.Range("D12").Value = '''14.87'
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = '''0.5902'
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").Value = '''4.389'
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").Value = '''75.30'
$ws.Range("E15").Value = '  +4.33%  '
$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '''1.000'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '25.562.93'
$ws.Range("E18").Value = '  +2.77%  '
$ws.Range("D19").Value = '''0.000006765'
$ws.Range("E19").Value = '  +3.19%  '
$ws.Range("D20").Value = '''11.47'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").Value = '1.884.99'
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").Value = '''4.440'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("E23").Value = '  +2.31%  '
$ws.Range("D24").Value = '''5.279'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").Value = '''136.73'
$ws.Range("E25").Value = '  +3.59%  '
$ws.Range("D26").Value = '''15.06'
$ws.Range("E26").Value = '  +2.24%  '
$ws.Range("D27").Value = '''1.391'
$ws.Range("D28").Value = '''1.723'
$ws.Range("E28").Value = '  +4.86%  '
$ws.Range("D29").Value = '''104.79'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("D30").Value = '''3.975'
$ws.Range("E30").Value = '  +7.30%  '
$ws.Range("D31").Value = '''0.07816'
$ws.Range("E31").Value = '  +1.05%  '
$ws.Range("D32").Value = '''3.654'
$ws.Range("E32").Value = '  +3.55%  '
$ws.Range("D33").Value = '''0.9990'
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").Value = '''0.04220'
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("D35").Value = '''2.620'
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("D36").Value = '''0.6100'
$ws.Range("E36").Value = '  +5.52%  '
$ws.Range("D37").Value = '''0.9533'
$ws.Range("E37").Value = '  +4.67%  '
$ws.Range("E38").Value = '  +3.36%  '
$ws.Range("D39").Value = '''0.8578'
$ws.Range("E39").Value = '  +4.11%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  +5.47%  '
$ws.Range("E42").Value = '  -5.06%  '
$ws.Range("D43").Value = '''96.47'
$ws.Range("E43").Value = '  -0.51%  '
$ws.Range("D44").Value = '''0.3768'
$ws.Range("E44").Value = '  +2.49%  '
$ws.Range("D45").Value = '''4.870'
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("D46").Value = '''0.1117'
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").Value = '''6.216'
$ws.Range("E47").Value = '  +2.92%  '
$ws.Range("D48").Value = '''0.05259'
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").Value = '''29.90'
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("D50").Value = '''7.379'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").Value = '''1.002'
$ws.Range("E51").Value = '  +0.11%  '
